# Update SVI to be rplthemes (overall measure); add table 4 to report.
# Now have comparison of no policy/policy and 75th percentile cutoff.
#
# The worksheet holds a "construct" -> "replace" lookup table (col A = internal
# variable name, col B = human readable label), used elsewhere for find/replace
# in a generated report. This edit:
#   - collapses the four standalone rpltheme1..4 rows into a single
#     "rplthemes" row labelled "SVI Overall Rank" (the SVI overall measure),
#   - fixes the stray "vaccination\\.50" construct name to "vaccination.50",
#   - appends a new "table 4" block of *.quarter.nonzero rows comparing the
#     no-policy/policy quarters against the 75th percentile cutoff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("construct", "replace"),
    @("changeinrate", "Change in COVID-19 case rate"),
    @("vaccination", "Vaccination"),
    @("etiquette", "Etiquette"),
    @("masks", "Masking"),
    @("physicaldistancing", "Physical distancing"),
    @("cohortingorstaggering", "Cohorting and/or staggering policy"),
    @("testingandscreening", "Testing and/or screening"),
    @("stayhome", "Stay home"),
    @("traceandquarantine", "Trace and quarantine"),
    @("cleaning", "Cleaning"),
    @("ventilation", "Ventilation"),
    @("percentamericanindianoralaskanative", "Percent American Indian/Alaska Native"),
    @("percentasian", "Percent Asian"),
    @("percentblackorafricanamerican", "Percent Black or African American"),
    @("percenthispaniclatino", "Percent Hispanic or Latino"),
    @("percentnativehawaiianorotherpacificislander", "Percent Native Hawaiian or other Pacific Islander"),
    @("percentnotspecified", "Percent no race specified"),
    @("percenttwoormoreraces", "Percent two or more races"),
    @("percentwhite", "Percent White"),
    @("percentstudentsfreereducedlunch", "Percent free and/or reduced lunch eligible"),
    @("rplthemes", "SVI Overall Rank"),
    @("vaccination.50", "50th percentile Vaccination"),
    @("etiquette.50", "50th percentile Etiquette"),
    @("masks.50", "50th percentile Masking"),
    @("physicaldistancing.50", "50th percentile Physical distancing"),
    @("cohortingorstaggering.50", "50th percentile Cohorting and/or staggering policy"),
    @("testingandscreening.50", "50th percentile Testing and/or screening"),
    @("stayhome.50", "50th percentile Stay home"),
    @("traceandquarantine.50", "50th percentile Trace and quarantine"),
    @("cleaning.50", "50th percentile Cleaning"),
    @("ventilation.50", "50th percentile Ventilation"),
    @("vaccination.75", "75th percentile Vaccination"),
    @("etiquette.75", "75th percentile Etiquette"),
    @("masks.75", "75th percentile Masking"),
    @("physicaldistancing.75", "75th percentile Physical distancing"),
    @("cohortingorstaggering.75", "75th percentile Cohorting and/or staggering"),
    @("testingandscreening.75", "75th percentile Testing and/or screening"),
    @("stayhome.75", "75th percentile Stay home"),
    @("traceandquarantine.75", "75th percentile Trace and quarantine"),
    @("cleaning.75", "75th percentile Cleaning"),
    @("ventilation.75", "75th percentile Ventilation"),
    @("vaccinationquarter.75", "75th percentile Vaccination"),
    @("etiquettequarter.75", "75th percentile Etiquette"),
    @("masksquarter.75", "75th percentile Masking"),
    @("physicaldistancingquarter.75", "75th percentile Physical distancing"),
    @("cohortingorstaggeringquarter.75", "75th percentile Cohorting and/or staggering policy"),
    @("testingandscreeningquarter.75", "75th percentile Testing and/or screening"),
    @("stayhomequarter.75", "75th percentile Stay home"),
    @("traceandquarantinequarter.75", "75th percentile Trace and quarantine"),
    @("cleaningquarter.75", "75th percentile Cleaning"),
    @("ventilationquarter.75", "75th percentile Ventilation"),
    @("vaccinationquarter", "Vaccination"),
    @("etiquettequarter", "Etiquette"),
    @("masksquarter", "Masking"),
    @("physicaldistancingquarter", "Physical distancing"),
    @("cohortingorstaggeringquarter", "Cohorting and/or staggering policy"),
    @("testingandscreeningquarter", "Testing and/or screening"),
    @("stayhomequarter", "Stay home"),
    @("traceandquarantinequarter", "Trace and quarantine"),
    @("cleaningquarter", "Cleaning"),
    @("ventilationquarter", "Ventilation"),
    @("vaccinationquarter.nonzero", "Vaccination"),
    @("etiquettequarter.nonzero", "Etiquette"),
    @("masksquarter.nonzero", "Masking"),
    @("physicaldistancingquarter.nonzero", "Physical distancing"),
    @("cohortingorstaggeringquarter.nonzero", "Cohorting and/or staggering policy"),
    @("testingandscreeningquarter.nonzero", "Testing and/or screening"),
    @("stayhomequarter.nonzero", "Stay home"),
    @("traceandquarantinequarter.nonzero", "Trace and quarantine"),
    @("cleaningquarter.nonzero", "Cleaning"),
    @("ventilationquarter.nonzero", "Ventilation")
)

$r = 1
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# Restore the view: scrolled so row 21 is at the top, with B23 selected
# (matches the sheetView/selection in the edited workbook).
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B23").Select()
